$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '72.271.49'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.06%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.645.64'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.17%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '590.79'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.83%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '175.08'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.12%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.643.45'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.11%  '
$ws.Range('E10').Value = '  -3.02%  '
$ws.Range('E11').Value = '  +1.45%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.356'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.64%  '
$ws.Range('E13').Value = '  -2.02%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.128.76'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.14%  '
$ws.Range('E15').Value = '  -2.67%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '72.167.74'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.01%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '25.99'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.31%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.681.04'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.39%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.32'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.04%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '8.00'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.91%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '369.74'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.29%  '
$ws.Range('E22').Value = '  -0.81%  '
$ws.Range('E23').Value = '  -0.48%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '71.39'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.32%  '
$ws.Range('E25').Value = '  -0.18%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.26'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E27').Value = '  -3.38%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.781.48'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.13%  '
$ws.Range('E29').Value = '  +0.05%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0954'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.40%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.05'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.15%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '499.19'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -5.16%  '
$ws.Range('E33').Value = '  -2.51%  '
$ws.Range('E34').Value = '  -0.76%  '
$ws.Range('E35').Value = '  -0.04%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '162.92'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.70%  '
$ws.Range('E37').Value = '  +3.12%  '
$ws.Range('E38').Value = '  +0.09%  '
$ws.Range('E39').Value = '  -0.95%  '
$ws.Range('E40').Value = '  -3.29%  '
$ws.Range('E41').Value = '  +0.02%  '
$ws.Range('E42').Value = '  -7.08%  '
$ws.Range('E43').Value = '  -3.40%  '
$ws.Range('E44').Value = '  -3.76%  '
$ws.Range('E46').Value = '  -0.45%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '153.87'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.61%  '
$ws.Range('E48').Value = '  +0.96%  '
$ws.Range('E49').Value = '  -1.26%  '
$ws.Range('E50').Value = '  -0.68%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0748'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.22%  '
